# Update the "COMPETENCES TECHNIQUES" skills list:
#  - remove the "Web : client" line entirely
#  - rename / retext the remaining skill lines
#  - append two brand-new skill lines at the end of the list
$d = $word.ActiveDocument

# 1) Remove the whole "Web : client" paragraph (its text AND its own
#    paragraph mark) without disturbing the empty section-break paragraph
#    that immediately precedes it. Including the paragraph-mark char (CR)
#    in the search string makes Find consume the mark atomically so the
#    previous (empty) paragraph is left completely untouched.
$cr = [char]13
$d.Content.Find.Execute("Web : client" + $cr, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 2) Langages : python, matlab, c, c++  ->  Langages : d, python, matlab, c, c++
$d.Content.Find.Execute("Langages : python, matlab, c, c++", $true, $false, $false, $false, $false, $true, 1, $false, "Langages : d, python, matlab, c, c++", 2) | Out-Null

# 3) MLOps : ...  ->  Data Science : ...
$d.Content.Find.Execute("MLOps : node.js, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit", $true, $false, $false, $false, $false, $true, 1, $false, "Data Science : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn", 2) | Out-Null

# 4) Visualisation : data engineering, tableau  ->  Visualisation : tableau
$d.Content.Find.Execute("Visualisation : data engineering, tableau", $true, $false, $false, $false, $false, $true, 1, $false, "Visualisation : tableau", 2) | Out-Null

# 5) ML/AI : ...  ->  Machine Learning : ...
$d.Content.Find.Execute("ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn", $true, $false, $false, $false, $false, $true, 1, $false, "Machine Learning : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit", 2) | Out-Null

# 6) Append two new skill lines ("Autres" and "Soft_Skills") right after
#    the "Machine Learning" line, matching its paragraph formatting.
$mlPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Machine Learning :") {
        $mlPara = $p
        break
    }
}

$mlPara.Range.InsertParagraphAfter()
$autresPara = $mlPara.Next()
$autresPara.Range.Text = "Autres :  technologies , si nécessaires"

$autresPara.Range.InsertParagraphAfter()
$softSkillsPara = $autresPara.Next()
$softSkillsPara.Range.Text = "Soft_Skills : initiative"
